# Updated cryptos list on Tue May  7 02:22:01 UTC 2024 with GitHub Actions
#
# This script updates the Price (column D) and Volume(1h) (column E) values
# for the crypto ranking sheet, and also reflects the re-ordering of three
# coin pairs (rows 13/14, 40/41, 43/44) whose rank positions swapped.
#
# Because many of the new "Price" strings look like plain numbers (e.g.
# "1.00", "592.56"), simply assigning them via .Value would let Excel
# auto-convert the cell to a Number type and silently reformat the text
# (e.g. "1.00" -> 1). To preserve the exact text representation used in
# the source data (including thousands-separator dots such as
# "63.731.84"), we force each Price cell to Text format before writing
# the value, then restore the cell style afterwards so no stray
# number-format / style metadata is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 13 / Row 14 swap: Avalanche moves up to rank 13 (row 13),
# ShibaInu moves down to rank 14 (row 14). Update Coin, Link, Price and
# Volume accordingly.
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D13" "37.53"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000241"
$ws.Range("E14").Value = "  -3.76%  "

# ---------------------------------------------------------------------
# Row 40 / Row 41 swap: Cosmos moves up to rank 40 (row 40),
# OKB moves down to rank 41 (row 41).
# ---------------------------------------------------------------------
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D40" "9.33"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D41" "50.96"
$ws.Range("E41").Value = "  -0.83%  "

# ---------------------------------------------------------------------
# Row 43 / Row 44 swap: Arweave moves up to rank 43 (row 43),
# TheGraph moves down to rank 44 (row 44).
# ---------------------------------------------------------------------
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D43" "42.05"
$ws.Range("E43").Value = "  +5.54%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D44" "0.289"
$ws.Range("E44").Value = "  -1.88%  "

# ---------------------------------------------------------------------
# Remaining rows: straightforward Price / Volume(1h) refreshes.
# ---------------------------------------------------------------------
Set-TextValue "D2" "63.731.84"
$ws.Range("E2").Value = "  -1.01%  "

Set-TextValue "D3" "3.091.54"
$ws.Range("E3").Value = "  -2.23%  "

$ws.Range("E4").Value = "  -0.22%  "

Set-TextValue "D5" "592.56"
$ws.Range("E5").Value = "  -0.10%  "

Set-TextValue "D6" "157.03"
$ws.Range("E6").Value = "  +6.70%  "

Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  -0.13%  "

Set-TextValue "D8" "0.543"
$ws.Range("E8").Value = "  +1.96%  "

Set-TextValue "D9" "3.091.51"
$ws.Range("E9").Value = "  -1.85%  "

$ws.Range("E10").Value = "  -3.99%  "

Set-TextValue "D11" "5.84"
$ws.Range("E11").Value = "  -1.11%  "

Set-TextValue "D12" "0.455"
$ws.Range("E12").Value = "  -0.81%  "

Set-TextValue "D15" "3.604.43"
$ws.Range("E15").Value = "  -1.96%  "

$ws.Range("E16").Value = "  -1.75%  "

Set-TextValue "D17" "63.768.79"
$ws.Range("E17").Value = "  -0.63%  "

Set-TextValue "D18" "7.16"
$ws.Range("E18").Value = "  -1.77%  "

Set-TextValue "D19" "3.089.97"
$ws.Range("E19").Value = "  -1.96%  "

Set-TextValue "D20" "479.42"
$ws.Range("E20").Value = "  +2.15%  "

Set-TextValue "D21" "14.62"
$ws.Range("E21").Value = "  +1.47%  "

Set-TextValue "D22" "0.712"
$ws.Range("E22").Value = "  -3.17%  "

Set-TextValue "D23" "7.58"
$ws.Range("E23").Value = "  -0.10%  "

Set-TextValue "D24" "2.41"
$ws.Range("E24").Value = "  +1.11%  "

Set-TextValue "D25" "81.38"
$ws.Range("E25").Value = "  -0.07%  "

Set-TextValue "D26" "12.87"
$ws.Range("E26").Value = "  -2.55%  "

Set-TextValue "D27" "10.24"
$ws.Range("E27").Value = "  +4.14%  "

Set-TextValue "D28" "1.00"
$ws.Range("E28").Value = "  -0.04%  "

Set-TextValue "D29" "7.48"
$ws.Range("E29").Value = "  +1.80%  "

Set-TextValue "D30" "2.68"
$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("E32").Value = "  -2.30%  "

Set-TextValue "D33" "0.114"
$ws.Range("E33").Value = "  -0.28%  "

Set-TextValue "D34" "27.33"
$ws.Range("E34").Value = "  -2.40%  "

Set-TextValue "D35" "0.0₃0855"
$ws.Range("E35").Value = "  -2.68%  "

Set-TextValue "D36" "3.52"
$ws.Range("E36").Value = "  +8.46%  "

$ws.Range("E37").Value = "  -1.14%  "

Set-TextValue "D38" "6.08"
$ws.Range("E38").Value = "  -1.88%  "

$ws.Range("E39").Value = "  -3.26%  "

Set-TextValue "D42" "447.71"
$ws.Range("E42").Value = "  -4.70%  "

Set-TextValue "D45" "0.0363"
$ws.Range("E45").Value = "  -3.26%  "

Set-TextValue "D46" "0.112"
$ws.Range("E46").Value = "  +3.87%  "

Set-TextValue "D47" "2.832.88"
$ws.Range("E47").Value = "  -2.70%  "

Set-TextValue "D48" "130.75"
$ws.Range("E48").Value = "  -2.41%  "

Set-TextValue "D49" "25.62"
$ws.Range("E49").Value = "  +5.92%  "

$ws.Range("E50").Value = "  +0.01%  "

Set-TextValue "D51" "2.27"
$ws.Range("E51").Value = "  +1.24%  "
